# Regenerate s_vals data to reflect filtered save-game data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-06-28)
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 6.048734245549538

# Row 3 (2023-03-17)
$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.034748368925986
